$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 (bold/centered/bordered like existing headers)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats

# New boolean columns F:H for rows 2-4
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = $false
